# B6-PowerPoint.pptx — commit "Fri, Mar 20, 2020  7:06:05 PM"
#
# The three tables in the deck (Component three slides) had their table
# style switched to a different built-in PowerPoint table style
# ({78DC9403-B70D-4487-AA55-C888F5A3FC4B}). Walk every slide, find any
# shape that carries a table, and re-apply the new style GUID — this
# mirrors picking a new style from the Table Styles gallery in the UI,
# which is surfaced on the object model as Table.ApplyStyle(StyleId).

$newTableStyleId = "{78DC9403-B70D-4487-AA55-C888F5A3FC4B}"

$p = $ppt.ActivePresentation

for ($slideIdx = 1; $slideIdx -le $p.Slides.Count; $slideIdx++) {
    $s = $p.Slides.Item($slideIdx)

    for ($shapeIdx = 1; $shapeIdx -le $s.Shapes.Count; $shapeIdx++) {
        $sh = $s.Shapes.Item($shapeIdx)

        if ($sh.HasTable) {
            $sh.Table.ApplyStyle($newTableStyleId)
        }
    }
}
